# Fixes during Regression Testing
# Update the "Date" (column B) result timestamps on the CCDeferredCC_27 and
# CMCAutopayCC_27 sheets to reflect the latest regression test run.

$wb = $excel.ActiveWorkbook

$wsCCDeferred = $wb.Worksheets.Item("CCDeferredCC_27")
$wsCCDeferred.Range("B2").Value = "Thu Nov 13 21:44:26 IST 2025"

$wsCMCAutopay = $wb.Worksheets.Item("CMCAutopayCC_27")
$wsCMCAutopay.Range("B2").Value = "Tue Nov 18 21:54:40 IST 2025"
